$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Date field: "8/" -> "2" + bookmark(_GoBack) + "/"   (08/11/2020 -> 02/11/2020)
# ------------------------------------------------------------------
$dateRng = $d.Content
$found = $dateRng.Find.Execute("08/11/2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $dateRng.Start
    # Replace the '8' (2nd character, right after the leading '0') with '2'
    $d.Range($start + 1, $start + 2).Text = "2"
    # Re-insert the document's "last edit" bookmark exactly between the new '2' and the following '/'
    $d.Bookmarks.Add("_GoBack", $d.Range($start + 2, $start + 2))
}

# ------------------------------------------------------------------
# 2) "Pseudocódigo" heading becomes bold
# ------------------------------------------------------------------
$boldRng = $d.Content
$foundBold = $boldRng.Find.Execute("Pseudocódigo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundBold) {
    $boldRng.Font.Bold = 1
    $boldRng.Paragraphs(1).Range.Font.Bold = 1
}

# ------------------------------------------------------------------
# 3) Merge the split "ESCRIB" / "IR" runs (and drop the stray bookmark
#    that used to sit between them) into a single continuous sentence.
# ------------------------------------------------------------------
$mergeRng = $d.Content
$mergeRng.Find.Execute(
    "ESCRIBIR. La lectura de datos se realiza, por defecto, desde el teclado, que es la entrada est" + [char]0x00E1 + "ndar del sistema. La escritura de datos se realiza, por defecto, en la pantalla, que es la salida est" + [char]0x00E1 + "ndar del sistema.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ESCRIBIR. La lectura de datos se realiza, por defecto, desde el teclado, que es la entrada est" + [char]0x00E1 + "ndar del sistema. La escritura de datos se realiza, por defecto, en la pantalla, que es la salida est" + [char]0x00E1 + "ndar del sistema.",
    2) | Out-Null
